$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$wb.Sheets.Item(1).Name = "Through 2022-09-27"

$ws.Range("A10").Value = "September (through 09-27)"

$ws.Range("B10").Value = 28
$ws.Range("C10").Value = 41
$ws.Range("D10").Value = 69
$ws.Range("E10").Value = 50
$ws.Range("F10").Value = 65
$ws.Range("G10").Value = 103
$ws.Range("H10").Value = 164
$ws.Range("I10").Value = 127

$ws.Range("B11").Value = 222
$ws.Range("C11").Value = 422
$ws.Range("D11").Value = 620
$ws.Range("E11").Value = 540
$ws.Range("F11").Value = 414
$ws.Range("G11").Value = 887
$ws.Range("H11").Value = 1234
$ws.Range("I11").Value = 1262
